$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.660.11"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.446.66"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.10"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.64"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.444.19"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").Value = "4.042.40"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.97"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "67.617.92"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "3.444.25"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.99"
$ws.Range("E20").Value = "  -6.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.99"
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.84"
$ws.Range("E22").Value = "  -3.63%  "
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E33").Value = "  -5.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.23"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("E37").Value = "  -7.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.88"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.882"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("E43").Value = "  -7.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.85"
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.99"
$ws.Range("E46").Value = "  -6.16%  "
$ws.Range("D47").Value = "2.695.82"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "323.10"
$ws.Range("E50").Value = "  -8.61%  "
$ws.Range("E51").Value = "  -4.93%  "
